$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.635.12'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '2.801.86'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'355.12"
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').Value = "'109.48"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = "'0.623"
$ws.Range('E9').Value = '  +5.75%  '
$ws.Range('D10').Value = "'40.18"
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = "'0.0838"
$ws.Range('E12').Value = '  -1.08%  '
$ws.Range('D13').Value = "'20.09"
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').Value = "'7.81"
$ws.Range('E14').Value = '  +3.16%  '
$ws.Range('D15').Value = '3.241.30'
$ws.Range('E15').Value = '  +0.49%  '
$ws.Range('D16').Value = '2.800.50'
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '51.604.55'
$ws.Range('E18').Value = '  -0.47%  '
$ws.Range('E19').Value = '  +3.44%  '
$ws.Range('E20').Value = '  +3.35%  '
$ws.Range('D21').Value = "'13.40"
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('D23').Value = "'70.43"
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = "'268.25"
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = "'2.77"
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = "'26.10"
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('E28').Value = '  +1.64%  '
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').Value = "'37.41"
$ws.Range('E30').Value = '  +7.87%  '
$ws.Range('D31').Value = "'2.23"
$ws.Range('E31').Value = '  +4.70%  '
$ws.Range('D32').Value = "'6.35"
$ws.Range('E32').Value = '  +10.62%  '
$ws.Range('D33').Value = "'52.19"
$ws.Range('E33').Value = '  +0.16%  '
$ws.Range('D34').Value = "'5.71"
$ws.Range('E34').Value = '  +10.02%  '
$ws.Range('D35').Value = "'0.0447"
$ws.Range('E35').Value = '  -4.93%  '
$ws.Range('D36').Value = "'0.0855"
$ws.Range('E36').Value = '  +1.30%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').Value = "'3.15"
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('E42').Value = '  -5.21%  '
$ws.Range('D43').Value = "'119.89"
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').Value = "'21.94"
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('D46').Value = '2.138.60'
$ws.Range('D47').Value = "'3.42"
$ws.Range('E47').Value = '  +4.79%  '
$ws.Range('E48').Value = '  +7.28%  '
$ws.Range('D49').Value = "'0.921"
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('E50').Value = '  +10.32%  '
$ws.Range('D51').Value = "'0.221"
$ws.Range('E51').Value = '  +16.30%  '
